$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.9950608814158476"
$ws.Range("E2").Value = [double]"0.9950608814158476"

$ws.Range("D3").Value = [double]"0.02451892750546465"
$ws.Range("E3").Value = [double]"0.02451892750546465"

$ws.Range("D4").Value = [double]"0.9990776111709648"
$ws.Range("E4").Value = [double]"0.9990776111709648"

$ws.Range("D5").Value = [double]"3.357927421125209E-18"
$ws.Range("E5").Value = [double]"3.357927421125209E-18"

$ws.Range("D6").Value = [double]"0.9999999996366451"
$ws.Range("E6").Value = [double]"0.9999999996366451"

$ws.Range("D7").Value = [double]"0.9999999999966329"
$ws.Range("E7").Value = [double]"3.367084389083175E-12"

$ws.Range("D8").Value = [double]"1.499720250130991E-07"
$ws.Range("E8").Value = [double]"0.999999850027975"

$ws.Range("D9").Value = [double]"0.0002771603847556075"
$ws.Range("E9").Value = [double]"0.9997228396152444"

$ws.Range("D10").Value = [double]"0.9999998913804018"
$ws.Range("E10").Value = [double]"1.0861959820474E-07"

$ws.Range("D11").Value = [double]"0.0003427357102554755"
$ws.Range("E11").Value = [double]"0.9996572642897446"
$ws.Range("F11").Value = [double]"6.594186305999756"
